# Added error pages, 1st draft of overview
# Updates the main data table: the "Start-End(Duration)" label text is
# corrected, and the financial columns (Measured Works .. Act. Revenue)
# are switched from pre-formatted text strings (e.g. "2 430 K") to plain
# numeric values (e.g. 2430), matching the new dashboard data source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duration label shared by all data rows (column C).
$ws.Cells.Replace("Apr 23 - Apr 24 (53)", "Apr 23-Apr 24 (53)")

# Replace the text-formatted "K" values in D:K with their numeric values.
$values = @{
    2 = @{ D = 2430;  E = 570;  F = 530;  G = 550;  H = 2000; I = 2000; J = 2000; K = 2000 }
    3 = @{ D = 12430; E = 1570; F = 1380; G = 1420; H = 3000; I = 3000; J = 3000; K = 3000 }
    4 = @{ D = 8500;  E = 1205; F = 1100; G = 1050; H = 4000; I = 4000; J = 4000; K = 4000 }
    5 = @{ D = 2850;  E = 280;  F = 250;  G = 245;  H = 5000; I = 5000; J = 5000; K = 5000 }
    6 = @{ D = 8850;  E = 1400; F = 1250; G = 1234; H = 6000; I = 6000; J = 6000; K = 6000 }
    7 = @{ D = 6800;  E = 985;  F = 900;  G = 905;  H = 7000; I = 7000; J = 7000; K = 7000 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}

# Reflect the author's last selection in the saved sheet view.
$ws.Range("L16").Select()
